$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even when it looks like a
# number (e.g. "6", "93", "42.85"), without touching cell styles/number
# formats. Direct .Value assignment on a purely-numeric-looking string makes
# Excel coerce the cell to a Number; going through a =TEXT() formula and then
# "flattening" it in place via Copy/PasteSpecial(xlPasteValues) keeps the
# text type while leaving styles.xml untouched.
function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '=TEXT("' + $escaped + '","@")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2: KS Bharat  (text, R=6, B=14, M=23, 4s=0, SR=42.85; 6s stays "0")
$ws.Range("A2").Value = "KS Bharat " + [char]0x2020
$ws.Range("B2").Value = "lbw b Dhawan"
Set-TextValue "C2" "6"
Set-TextValue "D2" "14"
Set-TextValue "E2" "23"
Set-TextValue "F2" "0"
Set-TextValue "H2" "42.85"

# Row 3: DB Prasanth (text, R=93, B=191, M=267, 4s=8, 6s=1, SR=48.69)
$ws.Range("A3").Value = "DB Prasanth"
$ws.Range("B3").Value = "not out"
Set-TextValue "C3" "93"
Set-TextValue "D3" "191"
Set-TextValue "E3" "267"
Set-TextValue "F3" "8"
Set-TextValue "G3" "1"
Set-TextValue "H3" "48.69"

# Row 4: GH Vihari (c) (text, R=68, B=149, M=243, 4s=9, SR=45.63; 6s stays "0")
$ws.Range("A4").Value = "GH Vihari (c)"
$ws.Range("B4").Value = "not out"
Set-TextValue "C4" "68"
Set-TextValue "D4" "149"
Set-TextValue "E4" "243"
Set-TextValue "F4" "9"
Set-TextValue "H4" "45.63"
